$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before the existing "GT1 Status" column (O), which
# shifts GT1..WT1 Status one column to the right (O->P, P->Q, Q->R, R->S, S->T)
$ws.Columns.Item(15).Insert()

# Header for the newly inserted column
$ws.Range("O1").Value = "Available Inertia [s]"

# Widen the new column to fit the longer header text (closest width the
# host's column-width pixel grid can represent to the authored 21.7109375)
$ws.Columns.Item(15).ColumnWidth = 20.83

# Fill the new column's data rows (2-12) with the available inertia value
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 15).Value = 9.6
}
